# Calibration edit for onshore wind
# Update the "Share of existing capacity" calibration row for "onshore wind es"
# (row 7) on the "CSC-CSCSoCECBiaSY" sheet from 0.2 to 0.33 for every year
# column (B:AE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$ws.Activate()

$ws.Range("B7:AE7").Value = 0.33

$ws.Range("B7:AE7").Select() | Out-Null
